# Replace the runs "elaborados " / "en la " / "tarde" (everything after
# "Documentos ") with a single run "modificados la próxima semana",
# while leaving the "Documentos " run (and the paragraph's own
# attributes: paraId/textId/rsid...) untouched.
#
# This engine's Find/Range text-editing collapses runs in a way that's
# hard to control directly when the edit touches the first couple of
# runs of a paragraph, so instead we: temporarily split the paragraph
# in two (so the original paragraph - with its original rsid/paraId
# attributes - becomes an isolated paragraph we can freely rewrite),
# rebuild the two pieces of text we want in each half, and then rejoin
# the two paragraphs by deleting the first paragraph's mark (which
# keeps the *second*, original paragraph's attributes/properties).

$d = $word.ActiveDocument

$prefix = "Documentos "
$newTail = "modificados la próxima semana"

# 1) Split the target paragraph into two: an empty paragraph in front of
#    it, and the (still fully intact, original) paragraph after it.
$d.Range(0, 0).InsertParagraphAfter()

$p1 = $d.Paragraphs(1)
$p2 = $d.Paragraphs(2)

# 2) Fill the brand-new leading paragraph with "Documentos ". It has a
#    single, fresh run, so this doesn't touch anything else.
$p1.Range.Text = $prefix

# 3) In the second paragraph (the original one, still carrying its
#    original rsid/paraId attributes and its original 4 runs), delete
#    the trailing runs one at a time working backwards from the end of
#    the paragraph. Deleting a run's exact full span down to nothing,
#    while it is the last run in the paragraph, removes it cleanly
#    without merging into/disturbing the runs before it.
$p2Start = $p2.Range.Start
$p2End = $p2.Range.End - 1   # exclude the paragraph mark

$d.Range($p2Start + $prefix.Length, $p2End).Text = ""

# 4) The second paragraph is now reduced to just its original first run
#    ("Documentos "), which is simultaneously the first and last run in
#    that paragraph. Overwrite it in place with the new tail text.
$d.Range($p2Start, $p2Start + $prefix.Length).Text = $newTail

# 5) Re-join the two paragraphs by deleting the leading paragraph's
#    mark. The surviving paragraph keeps the *second* paragraph's
#    properties (its original paraId/textId/rsid attributes), while the
#    text of both paragraphs is concatenated in order, giving:
#    "Documentos " + "modificados la próxima semana" as two runs in the
#    original paragraph.
$p1 = $d.Paragraphs(1)
$d.Range($p1.Range.End - 1, $p1.Range.End).Delete()

$d.Save()
